$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "tag_3" entry (row 4, column B) to "tag_333"
$ws.Range("B4").Value = "tag_333"

# Update the active selection to C9, as saved in the workbook
$ws.Range("C9").Select()
